# Weekly update: a new record (week) is added at the top of the data table
# (row 584), pushing every existing record down by one row. The former last
# row (702) ends up at row 703. Columns A, B, C, E, F, G, H, N, Q, R are
# constant across the whole data range, so only D, I, J, K, L, M, O, P vary
# and need to be set explicitly for the newly inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 584; this shifts rows 584:702 down to 585:703
# and carries the D-column date style (s="2") onto the new row automatically,
# and also grows the sheet's used range / dimension to A1:R703.
$ws.Rows("584:584").Insert()

# Populate the newly inserted row 584 with the new week's record.
$ws.Range("A584").Value = 9
$ws.Range("B584").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C584").Value = "Metropolitana"
$ws.Range("D584").Value = 45275
$ws.Range("E584").Value = 13
$ws.Range("F584").Value = 100112012
$ws.Range("G584").Value = "Espinaca"
$ws.Range("H584").Value = "Sin especificar"
$ws.Range("I584").Value = "Primera"
$ws.Range("J584").Value = 160
$ws.Range("K584").Value = 8000
$ws.Range("L584").Value = 10000
$ws.Range("M584").Value = 9000
$ws.Range("N584").Value = "$/cuna 10 kilos"
$ws.Range("O584").Value = "Provincia de Chacabuco"
$ws.Range("P584").Value = 900
$ws.Range("Q584").Value = 10
$ws.Range("R584").Value = "Hortaliza"
